$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A20 is a date-formatted-looking string ("11-10-2025") that must stay plain
# text (like the existing date cells above it), not get auto-converted into
# a date serial number. Enter it as a text formula, then paste-special just
# the value back onto itself so it is committed as a literal string without
# ever touching NumberFormat (which would otherwise leave a stray style).
$ws.Range("A20").Formula = '="11-10-2025"'
$ws.Range("A20").Copy()
$ws.Range("A20").PasteSpecial(-4163)

$ws.Range("B20").Value = "The price of gold in India today is ₹12,426 per gram for 24 karat gold, ₹11,390 per gram for 22 karat gold and ₹9,319 per gram for 18 karat gold (also called 999 gold)."
